$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Start of ramp / End of ramp triplet values for rows 3,4,5,7,8
$ws.Range("D3").Value = 1.1
$ws.Range("G3").Value = 0.5
$ws.Range("J3").Value = 3.5

$ws.Range("D4").Value = 1.7
$ws.Range("G4").Value = 9.2
$ws.Range("J4").Value = 7.2

$ws.Range("D5").Value = 10.8
$ws.Range("G5").Value = 3.3

$ws.Range("D7").Value = 5.2
$ws.Range("G7").Value = 3.3
$ws.Range("J7").Value = 3.4

$ws.Range("D8").Value = 1.3
$ws.Range("G8").Value = 3.9
$ws.Range("J8").Value = 3.1

# Update the active cell selection to I26
$ws.Range("I26").Select()
